$wb = $excel.ActiveWorkbook

# Sheet names: "Sheet2" tab is the active/visible one, "Sheet1" tab holds the
# dept_data/ind_data/bn1_data/singt_data rows being cleaned up.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- "Sheet1" tab: drop the now-unused reference-data labels in A31:A34 ---
# (this also causes the now-orphaned shared strings "dept_data"/"ind_data"/
# "bn1_data"/"singt_data" to be garbage collected on save, which re-packs the
# shared string table and shifts every other <v> index down accordingly.)
$ws1.Activate()
$ws1.Range("A31:A34").ClearContents()

# Scroll/position the frozen pane further down and leave the selection on A34
# (matches the new topLeftCell/selection recorded for this sheet's view).
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$ws1.Range("A34").Select()

# --- "Sheet2" tab: fill in the missing "x" marks for rows 23, 29, 30 (which
# previously had no marks at all) and rows 24-28 (missing only column D) ---
$ws2.Activate()

$ws2.Range("B23:G23").Value = "x"
$ws2.Range("F23:G23").WrapText = $true
$ws2.Rows.Item(23).RowHeight = 17

$ws2.Range("D24").Value = "x"
$ws2.Range("D25").Value = "x"
$ws2.Range("D26").Value = "x"
$ws2.Range("D27").Value = "x"
$ws2.Range("D28").Value = "x"

$ws2.Range("B29:G29").Value = "x"
$ws2.Range("F29:G29").WrapText = $true
$ws2.Rows.Item(29).RowHeight = 17

$ws2.Range("B30:G30").Value = "x"
$ws2.Range("F30:G30").WrapText = $true
$ws2.Rows.Item(30).RowHeight = 17

# Leave the selection on B11 for the visible/active sheet, with no frozen
# top-left scroll offset.
$ws2.Range("B11").Select()
